{"js": "// Replace the date heading and every \"NNN\u00f7N=\" drill cell with its updated\n// value. Each \"before\" string is unique in the document, so an exact,\n// case-sensitive search-and-replace is unambiguous for every entry.\nconst replacements = [\n  [\"2024-12-01 Sunday\", \"2024-12-02 Monday\"],\n  [\"347\u00f76=\", \"156\u00f75=\"],\n  [\"171\u00f77=\", \"629\u00f75=\"],\n  [\"821\u00f74=\", \"332\u00f77=\"],\n  [\"209\u00f73=\", \"325\u00f73=\"],\n  [\"276\u00f73=\", \"391\u00f76=\"],\n  [\"223\u00f77=\", \"374\u00f78=\"],\n  [\"499\u00f77=\", \"333\u00f76=\"],\n  [\"868\u00f74=\", \"173\u00f78=\"],\n  [\"543\u00f75=\", \"738\u00f77=\"],\n  [\"414\u00f73=\", \"931\u00f79=\"],\n  [\"841\u00f73=\", \"958\u00f75=\"],\n  [\"260\u00f77=\", \"118\u00f78=\"],\n  [\"492\u00f77=\", \"795\u00f73=\"],\n  [\"165\u00f77=\", \"349\u00f72=\"],\n  [\"698\u00f76=\", \"956\u00f78=\"],\n  [\"200\u00f75=\", \"174\u00f76=\"],\n  [\"879\u00f74=\", \"901\u00f79=\"],\n  [\"174\u00f72=\", \"690\u00f78=\"],\n  [\"548\u00f73=\", \"716\u00f77=\"],\n  [\"821\u00f72=\", \"826\u00f76=\"],\n  [\"915\u00f73=\", \"415\u00f79=\"],\n  [\"629\u00f79=\", \"817\u00f73=\"],\n  [\"742\u00f79=\", \"872\u00f79=\"],\n  [\"782\u00f77=\", \"917\u00f76=\"],\n  [\"539\u00f75=\", \"577\u00f74=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [find, replace] of replacements) {\n  const results = body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items/text\");\n  await context.sync();\n\n  for (const r of results.items) {\n    r.insertText(replace, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date heading and every \"NNN\u00f7N=\" drill cell with its updated\n# value. Each \"before\" string is unique in the document, so an exact,\n# case-sensitive search-and-replace is unambiguous for every entry.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"2024-12-01 Sunday\", \"2024-12-02 Monday\"),\n  @(\"347\u00f76=\", \"156\u00f75=\"),\n  @(\"171\u00f77=\", \"629\u00f75=\"),\n  @(\"821\u00f74=\", \"332\u00f77=\"),\n  @(\"209\u00f73=\", \"325\u00f73=\"),\n  @(\"276\u00f73=\", \"391\u00f76=\"),\n  @(\"223\u00f77=\", \"374\u00f78=\"),\n  @(\"499\u00f77=\", \"333\u00f76=\"),\n  @(\"868\u00f74=\", \"173\u00f78=\"),\n  @(\"543\u00f75=\", \"738\u00f77=\"),\n  @(\"414\u00f73=\", \"931\u00f79=\"),\n  @(\"841\u00f73=\", \"958\u00f75=\"),\n  @(\"260\u00f77=\", \"118\u00f78=\"),\n  @(\"492\u00f77=\", \"795\u00f73=\"),\n  @(\"165\u00f77=\", \"349\u00f72=\"),\n  @(\"698\u00f76=\", \"956\u00f78=\"),\n  @(\"200\u00f75=\", \"174\u00f76=\"),\n  @(\"879\u00f74=\", \"901\u00f79=\"),\n  @(\"174\u00f72=\", \"690\u00f78=\"),\n  @(\"548\u00f73=\", \"716\u00f77=\"),\n  @(\"821\u00f72=\", \"826\u00f76=\"),\n  @(\"915\u00f73=\", \"415\u00f79=\"),\n  @(\"629\u00f79=\", \"817\u00f73=\"),\n  @(\"742\u00f79=\", \"872\u00f79=\"),\n  @(\"782\u00f77=\", \"917\u00f76=\"),\n  @(\"539\u00f75=\", \"577\u00f74=\")\n)\n\nforeach ($pair in $pairs) {\n  $findText = $pair[0]\n  $replaceText = $pair[1]\n\n  # Find.Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards,\n  #   MatchSoundsLike, MatchAllWordForms, Forward, Wrap, Format,\n  #   ReplaceWith, Replace)\n  # Wrap:=1 -> wdFindContinue, Replace:=2 -> wdReplaceAll\n  $rng = $d.Content\n  $rng.Find.ClearFormatting()\n  $rng.Find.Replacement.ClearFormatting()\n  $rng.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n}\n"}
